$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# The "Ready for handoff" status string is shared across the Overview
# sheet (columns B & C) and the Status column (C) on both the zh-cn and
# de-de detail sheets for the 675e92ad... row (row 3). Update all of them
# so every cell that used to read "Ready for handoff" now reads
# "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus

# zh-cn sheet: add Error Detail (column K) for row 3.
$wsZhCn.Range("K3").Value = "Handback file name: vqej1znn.kyj is different with handoff file name: 675e92ad-4479-49f1-b054-9903a75dbf2a.17a1217f5c920d608f83edb141967d037ff78ab4.zh-cn."

# de-de sheet: add Error Detail (column K) for row 3.
$wsDeDe.Range("K3").Value = "Handback file name: vqej1znn.kyj is different with handoff file name: 675e92ad-4479-49f1-b054-9903a75dbf2a.17a1217f5c920d608f83edb141967d037ff78ab4.de-de."
